# Updates loading_percent values on Sheet1 (rows 2-25, columns C-G and I-L)
# to reflect the re-run results for the "case with 380 kV" scenario.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Each inner array: row, C, D, E, F, G, I, J, K, L
$data = @(
    @(2,  6.122469343576001, 4.180267855248101, 12.10055612961523, 58.88169246593258, 3.81370664019963,  44.89219691801183, 11.33704543409912, 26.74595557374489, 10.59323797635611),
    @(3,  6.104175000773907, 4.168658371637316, 12.1109545550067,  58.92316984943575, 3.818075198294558, 44.94553339962354, 11.36319620910437, 26.50691985857935, 10.61470512249857),
    @(4,  6.09367375611464,  4.16160006536476,  12.11897750807961, 58.96310536059362, 3.82089436713839,  44.98878250786377, 11.38061867545793, 26.36553396721238, 10.62904323288335),
    @(5,  6.08958163469481,  4.158742609013119, 12.12265945807638, 58.98300199596405, 3.822077754193025, 45.00903554344015, 11.38806209132759, 26.30932809107871, 10.63517748018583),
    @(6,  6.088913541767963, 4.15826930691388,  12.12329577452807, 58.98652412459798, 3.822276345702697, 45.01255694917364, 11.38931882399145, 26.30008190605362, 10.63621367485425),
    @(7,  6.093617806041173, 4.161561450638565, 12.11902549310526, 58.96335904802304, 3.820910186625404, 44.98904502042191, 11.38071766843733, 26.36477017355172, 10.62912478134374),
    @(8,  6.116011613239538, 4.176251140797118, 12.10380184152336, 58.8929821167393,  3.815184598455303, 44.90840122837378, 11.34577886183814, 26.66245345800626, 10.60039989140593),
    @(9,  6.165582507670485, 4.20555828513172,  12.08692101381932, 58.87039825172492, 3.805036277750079, 44.83406032978602, 11.2880930261111,  27.28632419827622, 10.5532363466677),
    @(10, 6.205246512065379, 4.2273243840406,   12.08239017987817, 58.92492484198056, 3.798229519926387, 44.83113868228679, 11.25230262861282, 27.7655202403142,  10.52415170388368),
    @(11, 6.223949562169937, 4.23726231481748,  12.08202909043969, 58.96529574587767, 3.795272013150766, 44.8411437976658,  11.23744981341041, 27.98724480825391, 10.5121245489294),
    @(12, 6.231122768595449, 4.241029522463126, 12.08213589815161, 58.98282813954114, 3.794171911805047, 44.84656923462067, 11.23203069386163, 28.07167512739879, 10.50774289821442),
    @(13, 6.229573919418572, 4.240218032856,    12.0821020788352,  58.97895230677776, 3.794407957959762, 44.84532787932199, 11.23318866792595, 28.05347183455362, 10.50867888585247),
    @(14, 6.224537910812717, 4.237572173963823, 12.08203300207822, 58.96669312799072, 3.795181110252409, 44.84155731421928, 11.23699986476355, 27.99418195811005, 10.51176060702233),
    @(15, 6.221464907775164, 4.235951982046175, 12.08202237886026, 58.9594765212567,  3.795657268419493, 44.83946106692684, 11.23936106899338, 27.95792409276973, 10.51367074120423),
    @(16, 6.204037150231121, 4.226675537660381, 12.08244790990897, 58.92260033693082, 3.798425583905003, 44.83071351180612, 11.2533020310698,  27.75109928598964, 10.52496190019363),
    @(17, 6.193511861228632, 4.220993133226798, 12.0831437446772,  58.9039701211044,  3.800159341151098, 44.82825539102582, 11.2622201501156,  27.62512668954966, 10.53219671394379),
    @(18, 6.187520337640097, 4.217728270312873, 12.08370404272795, 58.8947194538265,  3.801169636436921, 44.82790803549138, 11.26748407062819, 27.55302639313048, 10.53647129369024),
    @(19, 6.185502541850372, 4.216623486048238, 12.08392126927092, 58.89183871396057, 3.801513956345297, 44.82797336875859, 11.2692894389352,  27.52867768547923, 10.53793806263944),
    @(20, 6.194625873230758, 4.221597681559722, 12.08305311140504, 58.90580167106546, 3.799973426566598, 44.82840663095791, 11.26125688623588, 27.63850036299001, 10.53141483110818),
    @(21, 6.226014678257673, 4.238349230592055, 12.08204668952036, 58.97023298192025, 3.794953479238393, 44.84262035161827, 11.23587485223335, 28.01158468178213, 10.51085074425375),
    @(22, 6.247056183890528, 4.249319639628172, 12.08280802251081, 59.02542853098392, 3.791788247267879, 44.86145222921675, 11.22048294493385, 28.25811280974521, 10.4984178303863),
    @(23, 6.235779079294132, 4.24346291399717,  12.0822721590594,  58.99477068845292, 3.793467058652782, 44.85052620717773, 11.22858842333816, 28.12631230627234, 10.50496147885105),
    @(24, 6.194122042438396, 4.22132435884555,  12.08309358748543, 58.90496907997389, 3.800057436380978, 44.82833493603162, 11.26169195180519, 27.63245311382621, 10.53176796152722),
    @(25, 6.151587322977268, 4.197583694429017, 12.090102457324,   58.86406559098685, 3.80766701675677,  44.84513293839549, 11.30254032045498, 27.11362304759333, 10.56501627319551)
)

# Column letters in the same order as the values stored in each row's array
# (note: column H is intentionally skipped - it stays 0 and is not touched)
$cols = @("C", "D", "E", "F", "G", "I", "J", "K", "L")

foreach ($rowData in $data) {
    $row = $rowData[0]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $col = $cols[$i]
        $value = $rowData[$i + 1]
        $ws.Range("$col$row").Value = $value
    }
}
